# Add two new columns "I0" (col I) and "IF" (col J) after the existing
# "IP" column (col H), mirroring its header style, and fill in the data
# rows with the recorded values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the bold/centered/bordered header style used by the other
# header cells (copy formatting from H1, which already carries it).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (2-37) ---------------------------------------------------
$arr = New-Object 'object[,]' 36,2
$arr[0,0] = 1
$arr[0,1] = 5
$arr[1,0] = 1
$arr[1,1] = 5
$arr[2,0] = 1
$arr[2,1] = 5
$arr[3,0] = 1
$arr[3,1] = 5
$arr[4,0] = 1
$arr[4,1] = 6
$arr[5,0] = 1
$arr[5,1] = 6
$arr[6,0] = 1
$arr[6,1] = 5
$arr[7,0] = 1
$arr[7,1] = 5
$arr[8,0] = 1
$arr[8,1] = 4
$arr[9,0] = 2
$arr[9,1] = 5
$arr[10,0] = 1
$arr[10,1] = 5
$arr[11,0] = 1
$arr[11,1] = 5
$arr[12,0] = 1
$arr[12,1] = 6
$arr[13,0] = 1
$arr[13,1] = 4
$arr[14,0] = 1
$arr[14,1] = 5
$arr[15,0] = 1
$arr[15,1] = 4
$arr[16,0] = 1
$arr[16,1] = 6
$arr[17,0] = 1
$arr[17,1] = 5
$arr[18,0] = 1
$arr[18,1] = 6
$arr[19,0] = 1
$arr[19,1] = 5
$arr[20,0] = 1
$arr[20,1] = 5
$arr[21,0] = 1
$arr[21,1] = 5
$arr[22,0] = 1
$arr[22,1] = 8
$arr[23,0] = 1
$arr[23,1] = 6
$arr[24,0] = 1
$arr[24,1] = 7
$arr[25,0] = 1
$arr[25,1] = 6
$arr[26,0] = 1
$arr[26,1] = 3
$arr[27,0] = 1
$arr[27,1] = 6
$arr[28,0] = 1
$arr[28,1] = 3
$arr[29,0] = 1
$arr[29,1] = 5
$arr[30,0] = 1
$arr[30,1] = 6
$arr[31,0] = 1
$arr[31,1] = 4
$arr[32,0] = 1
$arr[32,1] = 5
$arr[33,0] = 1
$arr[33,1] = 3
$arr[34,0] = 1
$arr[34,1] = 3
$arr[35,0] = 1
$arr[35,1] = 1

$ws.Range("I2:J37").Value = $arr
